# ---------------------------------------------------------------
# Sheet "存款" (deposits, sheet index 3): add bank / deposit_type /
# currency columns + shift totals + append metadata columns
# (property_category, category, date, legislator_name,
#  legislator_id, source_file, index).
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

function Set-DepositHeader($col, $text) {
    $c = $ws.Cells.Item(1, $col)
    $c.Value = $text
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# Row 1 already has B1=bank, C1=deposit_type, D1=currency, E1=owner,
# F1=total headers from the original sheet; only G1..M1 are new.
Set-DepositHeader 7 "property_category"
Set-DepositHeader 8 "category"
Set-DepositHeader 9 "date"
Set-DepositHeader 10 "legislator_name"
Set-DepositHeader 11 "legislator_id"
Set-DepositHeader 12 "source_file"
Set-DepositHeader 13 "index"

# Row 2 (index 48): 臺灣銀行永和分行 / 活期存款 / 新臺幣
$ws.Cells.Item(2, 2).Value = "臺灣銀行永和分行"
$ws.Cells.Item(2, 3).Value = "活期存款"
$ws.Cells.Item(2, 4).Value = "新臺幣"
$ws.Cells.Item(2, 6).Value = 4161600
$ws.Cells.Item(2, 7).Value = "deposit"
$ws.Cells.Item(2, 8).Value = "normal"
$ws.Cells.Item(2, 9).Value = "2012-04-05"
$ws.Cells.Item(2, 10).Value = "陳鎮湘"
$ws.Cells.Item(2, 11).Value = 1754
$ws.Cells.Item(2, 12).Value = "tmpc08e1"
$ws.Cells.Item(2, 13).Value = 48

# Row 3 (index 49): 臺灣銀行龍山分行 / 活期存款 / 新臺幣
$ws.Cells.Item(3, 2).Value = "臺灣銀行龍山分行"
$ws.Cells.Item(3, 3).Value = "活期存款"
$ws.Cells.Item(3, 4).Value = "新臺幣"
$ws.Cells.Item(3, 6).Value = 461344
$ws.Cells.Item(3, 7).Value = "deposit"
$ws.Cells.Item(3, 8).Value = "normal"
$ws.Cells.Item(3, 9).Value = "2012-04-05"
$ws.Cells.Item(3, 10).Value = "陳鎮湘"
$ws.Cells.Item(3, 11).Value = 1754
$ws.Cells.Item(3, 12).Value = "tmpc08e1"
$ws.Cells.Item(3, 13).Value = 49

# Row 4 (index 50): 臺灣銀行永和分行 / 公教優惠儲蓄存款 / 新臺幣
$ws.Cells.Item(4, 2).Value = "臺灣銀行永和分行"
$ws.Cells.Item(4, 3).Value = "公教優惠儲蓄存款"
$ws.Cells.Item(4, 4).Value = "新臺幣"
$ws.Cells.Item(4, 6).Value = 522300
$ws.Cells.Item(4, 7).Value = "deposit"
$ws.Cells.Item(4, 8).Value = "normal"
$ws.Cells.Item(4, 9).Value = "2012-04-05"
$ws.Cells.Item(4, 10).Value = "陳鎮湘"
$ws.Cells.Item(4, 11).Value = 1754
$ws.Cells.Item(4, 12).Value = "tmpc08e1"
$ws.Cells.Item(4, 13).Value = 50

# Row 5 (index 51): 中華郵政股份有限公司中和宜安郵局 / 活期存款 / 新臺幣
$ws.Cells.Item(5, 2).Value = "中華郵政股份有限公司中和宜安郵局"
$ws.Cells.Item(5, 3).Value = "活期存款"
$ws.Cells.Item(5, 4).Value = "新臺幣"
$ws.Cells.Item(5, 6).Value = 358750
$ws.Cells.Item(5, 7).Value = "deposit"
$ws.Cells.Item(5, 8).Value = "normal"
$ws.Cells.Item(5, 9).Value = "2012-04-05"
$ws.Cells.Item(5, 10).Value = "陳鎮湘"
$ws.Cells.Item(5, 11).Value = 1754
$ws.Cells.Item(5, 12).Value = "tmpc08e1"
$ws.Cells.Item(5, 13).Value = 51

# Row 6 (index 52): 國防部主計局同袍儲蓄會 / 定期存款 / 新臺幣
$ws.Cells.Item(6, 2).Value = "國防部主計局同袍儲蓄會"
$ws.Cells.Item(6, 3).Value = "定期存款"
$ws.Cells.Item(6, 4).Value = "新臺幣"
$ws.Cells.Item(6, 6).Value = 80000
$ws.Cells.Item(6, 7).Value = "deposit"
$ws.Cells.Item(6, 8).Value = "normal"
$ws.Cells.Item(6, 9).Value = "2012-04-05"
$ws.Cells.Item(6, 10).Value = "陳鎮湘"
$ws.Cells.Item(6, 11).Value = 1754
$ws.Cells.Item(6, 12).Value = "tmpc08e1"
$ws.Cells.Item(6, 13).Value = 52

# Row 7 (index 53): 國防部主計局同袍儲蓄會 / 定期存款 / 新臺幣
$ws.Cells.Item(7, 2).Value = "國防部主計局同袍儲蓄會"
$ws.Cells.Item(7, 3).Value = "定期存款"
$ws.Cells.Item(7, 4).Value = "新臺幣"
$ws.Cells.Item(7, 6).Value = 550000
$ws.Cells.Item(7, 7).Value = "deposit"
$ws.Cells.Item(7, 8).Value = "normal"
$ws.Cells.Item(7, 9).Value = "2012-04-05"
$ws.Cells.Item(7, 10).Value = "陳鎮湘"
$ws.Cells.Item(7, 11).Value = 1754
$ws.Cells.Item(7, 12).Value = "tmpc08e1"
$ws.Cells.Item(7, 13).Value = 53

# Row 8 (index 54): 中華郵政股份有限公司中和宜安郵局 / 活期存款 / 新臺幣
$ws.Cells.Item(8, 2).Value = "中華郵政股份有限公司中和宜安郵局"
$ws.Cells.Item(8, 3).Value = "活期存款"
$ws.Cells.Item(8, 4).Value = "新臺幣"
$ws.Cells.Item(8, 6).Value = 543375
$ws.Cells.Item(8, 7).Value = "deposit"
$ws.Cells.Item(8, 8).Value = "normal"
$ws.Cells.Item(8, 9).Value = "2012-04-05"
$ws.Cells.Item(8, 10).Value = "陳鎮湘"
$ws.Cells.Item(8, 11).Value = 1754
$ws.Cells.Item(8, 12).Value = "tmpc08e1"
$ws.Cells.Item(8, 13).Value = 54

# Row 9 (index 55): 萬泰商業銀行 / 活期存款 / 新臺幣
$ws.Cells.Item(9, 2).Value = "萬泰商業銀行"
$ws.Cells.Item(9, 3).Value = "活期存款"
$ws.Cells.Item(9, 4).Value = "新臺幣"
$ws.Cells.Item(9, 6).Value = 91991
$ws.Cells.Item(9, 7).Value = "deposit"
$ws.Cells.Item(9, 8).Value = "normal"
$ws.Cells.Item(9, 9).Value = "2012-04-05"
$ws.Cells.Item(9, 10).Value = "陳鎮湘"
$ws.Cells.Item(9, 11).Value = 1754
$ws.Cells.Item(9, 12).Value = "tmpc08e1"
$ws.Cells.Item(9, 13).Value = 55

# Row 10 (index 56): 第一商業銀行 / 活期存款 / 新臺幣
$ws.Cells.Item(10, 2).Value = "第一商業銀行"
$ws.Cells.Item(10, 3).Value = "活期存款"
$ws.Cells.Item(10, 4).Value = "新臺幣"
$ws.Cells.Item(10, 6).Value = 78533
$ws.Cells.Item(10, 7).Value = "deposit"
$ws.Cells.Item(10, 8).Value = "normal"
$ws.Cells.Item(10, 9).Value = "2012-04-05"
$ws.Cells.Item(10, 10).Value = "陳鎮湘"
$ws.Cells.Item(10, 11).Value = 1754
$ws.Cells.Item(10, 12).Value = "tmpc08e1"
$ws.Cells.Item(10, 13).Value = 56

# Row 11 (index 57): 國泰世華商業銀行 / 活期存款 / 新臺幣
$ws.Cells.Item(11, 2).Value = "國泰世華商業銀行"
$ws.Cells.Item(11, 3).Value = "活期存款"
$ws.Cells.Item(11, 4).Value = "新臺幣"
$ws.Cells.Item(11, 6).Value = 10000
$ws.Cells.Item(11, 7).Value = "deposit"
$ws.Cells.Item(11, 8).Value = "normal"
$ws.Cells.Item(11, 9).Value = "2012-04-05"
$ws.Cells.Item(11, 10).Value = "陳鎮湘"
$ws.Cells.Item(11, 11).Value = 1754
$ws.Cells.Item(11, 12).Value = "tmpc08e1"
$ws.Cells.Item(11, 13).Value = 57

# Row 12 (index 58): 渣打國際商業銀行 / 活期存款 / 新臺幣
$ws.Cells.Item(12, 2).Value = "渣打國際商業銀行"
$ws.Cells.Item(12, 3).Value = "活期存款"
$ws.Cells.Item(12, 4).Value = "新臺幣"
$ws.Cells.Item(12, 6).Value = 112695.89
$ws.Cells.Item(12, 7).Value = "deposit"
$ws.Cells.Item(12, 8).Value = "normal"
$ws.Cells.Item(12, 9).Value = "2012-04-05"
$ws.Cells.Item(12, 10).Value = "陳鎮湘"
$ws.Cells.Item(12, 11).Value = 1754
$ws.Cells.Item(12, 12).Value = "tmpc08e1"
$ws.Cells.Item(12, 13).Value = 58

# Row 13 (index 59): 渣打國際商業銀行 / 活期存款 / 美金
$ws.Cells.Item(13, 2).Value = "渣打國際商業銀行"
$ws.Cells.Item(13, 3).Value = "活期存款"
$ws.Cells.Item(13, 4).Value = "美金"
$ws.Cells.Item(13, 6).Value = 1049651.12
$ws.Cells.Item(13, 7).Value = "deposit"
$ws.Cells.Item(13, 8).Value = "normal"
$ws.Cells.Item(13, 9).Value = "2012-04-05"
$ws.Cells.Item(13, 10).Value = "陳鎮湘"
$ws.Cells.Item(13, 11).Value = 1754
$ws.Cells.Item(13, 12).Value = "tmpc08e1"
$ws.Cells.Item(13, 13).Value = 59

# Row 14 (index 60): 渣打國際商業銀行 / 定期存款 / 美金
$ws.Cells.Item(14, 2).Value = "渣打國際商業銀行"
$ws.Cells.Item(14, 3).Value = "定期存款"
$ws.Cells.Item(14, 4).Value = "美金"
$ws.Cells.Item(14, 6).Value = 294900
$ws.Cells.Item(14, 7).Value = "deposit"
$ws.Cells.Item(14, 8).Value = "normal"
$ws.Cells.Item(14, 9).Value = "2012-04-05"
$ws.Cells.Item(14, 10).Value = "陳鎮湘"
$ws.Cells.Item(14, 11).Value = 1754
$ws.Cells.Item(14, 12).Value = "tmpc08e1"
$ws.Cells.Item(14, 13).Value = 60

# Row 15 (index 61): 兆豐國際商業銀行永和分行 / 活期存款 / 美金
$ws.Cells.Item(15, 2).Value = "兆豐國際商業銀行永和分行"
$ws.Cells.Item(15, 3).Value = "活期存款"
$ws.Cells.Item(15, 4).Value = "美金"
$ws.Cells.Item(15, 6).Value = 182838
$ws.Cells.Item(15, 7).Value = "deposit"
$ws.Cells.Item(15, 8).Value = "normal"
$ws.Cells.Item(15, 9).Value = "2012-04-05"
$ws.Cells.Item(15, 10).Value = "陳鎮湘"
$ws.Cells.Item(15, 11).Value = 1754
$ws.Cells.Item(15, 12).Value = "tmpc08e1"
$ws.Cells.Item(15, 13).Value = 61

# Row 16 (index 62): 兆豐國際商業銀行南台北分行 / 定期存款 / 歐元
$ws.Cells.Item(16, 2).Value = "兆豐國際商業銀行南台北分行"
$ws.Cells.Item(16, 3).Value = "定期存款"
$ws.Cells.Item(16, 4).Value = "歐元"
$ws.Cells.Item(16, 6).Value = 440000
$ws.Cells.Item(16, 7).Value = "deposit"
$ws.Cells.Item(16, 8).Value = "normal"
$ws.Cells.Item(16, 9).Value = "2012-04-05"
$ws.Cells.Item(16, 10).Value = "陳鎮湘"
$ws.Cells.Item(16, 11).Value = 1754
$ws.Cells.Item(16, 12).Value = "tmpc08e1"
$ws.Cells.Item(16, 13).Value = 62
